$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B26").Value = 6488
$ws.Range("D26").Value = 6046810
$ws.Range("E26").Value = 931.999075215783
$ws.Range("F26").Value = 9.650160554334963
$ws.Range("H26").Value = 25.92502196520428
